$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("N21").ClearContents()
$ws.Range("H23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("N23").ClearContents()
$ws.Range("H34").Value = 2377.7144
$ws.Range("I34").Value = 2377.7144
$ws.Range("K34").Value = 2377.7144
$ws.Range("M34").Value = -2174.7144
$ws.Range("H36").Value = 2377.7144
$ws.Range("I36").Value = 2377.7144
$ws.Range("K36").Value = 2377.7144
$ws.Range("M36").Value = -1662.7144
$ws.Range("H44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").ClearContents()
$ws.Range("H47").Value = 14500
$ws.Range("I47").Value = 9000
$ws.Range("K47").Value = 9000
$ws.Range("M47").Value = -8028
$ws.Range("H74").Value = 3550.8333
$ws.Range("I74").Value = 3166.611
$ws.Range("J74").Value = 3935.0557
$ws.Range("K74").Value = 3166.611
$ws.Range("L74").Value = 3935.0557
$ws.Range("M74").Value = -2230.611
$ws.Range("N74").Value = -5807.0557
$ws.Range("H77").Value = 3550.8333
$ws.Range("I77").Value = 3166.611
$ws.Range("J77").Value = 3935.0557
$ws.Range("K77").Value = 15833.055
$ws.Range("L77").Value = 19675.2785
$ws.Range("M77").Value = -11153.055
$ws.Range("N77").Value = -29035.2785
$ws.Range("H96").Value = 876.4545000000001
$ws.Range("I96").Value = 340.2857
$ws.Range("J96").Value = 1814.75
$ws.Range("K96").Value = 1020.8571
$ws.Range("L96").Value = 5444.25
$ws.Range("M96").Value = 352.1428999999999
$ws.Range("N96").Value = -8190.25
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1310.2941
$ws.Range("I2").Value = 846.8333
$ws.Range("J2").Value = 2422.6
$ws.Range("K2").Value = 846.8333
$ws.Range("L2").Value = 2422.6
$ws.Range("M2").Value = -733.8333
$ws.Range("N2").Value = -2648.6
$ws.Range("H32").Value = 8337.366
$ws.Range("I32").Value = 7315.52
$ws.Range("K32").Value = 7315.52
$ws.Range("M32").Value = -7028.52
$ws.Range("H44").Value = 28124.625
$ws.Range("J44").Value = 28124.625
$ws.Range("L44").Value = 28124.625
$ws.Range("N44").Value = -29100.625
$ws.Range("H47").Value = 24041
$ws.Range("J47").Value = 24041
$ws.Range("L47").Value = 24041
$ws.Range("N47").Value = -25491
$ws.Range("H53").Value = 19526.2
$ws.Range("I53").Value = 4772.5
$ws.Range("J53").Value = 29362
$ws.Range("K53").Value = 4772.5
$ws.Range("L53").Value = 29362
$ws.Range("M53").Value = -4090.5
$ws.Range("N53").Value = -30726
$ws.Range("H116").Value = 1310.2941
$ws.Range("I116").Value = 846.8333
$ws.Range("J116").Value = 2422.6
$ws.Range("K116").Value = 846.8333
$ws.Range("L116").Value = 2422.6
$ws.Range("M116").Value = 1447.1667
$ws.Range("N116").Value = -7010.6
$ws.Range("H132").Value = 20837234
$ws.Range("I132").Value = 35718256
$ws.Range("J132").Value = 3802.4
$ws.Range("K132").Value = 107154768
$ws.Range("L132").Value = 11407.2
$ws.Range("M132").Value = -107152238
$ws.Range("N132").Value = -16467.2
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1310.2941
$ws.Range("I3").Value = 846.8333
$ws.Range("J3").Value = 2422.6
$ws.Range("K3").Value = 846.8333
$ws.Range("L3").Value = 2422.6
$ws.Range("M3").Value = -732.8333
$ws.Range("N3").Value = -2650.6
$ws.Range("H94").Value = 782.7083
$ws.Range("I94").Value = 655.4761999999999
$ws.Range("J94").Value = 1673.3334
$ws.Range("K94").Value = 655.4761999999999
$ws.Range("L94").Value = 1673.3334
$ws.Range("M94").Value = -204.4761999999999
$ws.Range("N94").Value = -2575.3334
$ws.Range("H97").Value = 13106.75
$ws.Range("I97").Value = 6214
$ws.Range("J97").Value = 19999.5
$ws.Range("K97").Value = 6214
$ws.Range("L97").Value = 19999.5
$ws.Range("M97").Value = -5223
$ws.Range("N97").Value = -21981.5
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 2361
$ws.Range("I99").Value = 1476.6666
$ws.Range("K99").Value = 1476.6666
$ws.Range("M99").Value = 21.33339999999998
$ws.Range("H107").Value = 805.7143
$ws.Range("I107").Value = 890
$ws.Range("J107").Value = 300
$ws.Range("K107").Value = 890
$ws.Range("L107").Value = 300
$ws.Range("M107").Value = 1030
$ws.Range("N107").Value = -4140
$ws.Range("H116").Value = 44000
$ws.Range("J116").Value = 44000
$ws.Range("L116").Value = 44000
$ws.Range("N116").Value = -53178
$ws.Range("H126").Value = 2361
$ws.Range("I126").Value = 1476.6666
$ws.Range("K126").Value = 4429.9998
$ws.Range("M126").Value = -1959.9998
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1530.8846
$ws.Range("I5").Value = 321.14285
$ws.Range("J5").Value = 2942.25
$ws.Range("K5").Value = 963.4285500000001
$ws.Range("L5").Value = 8826.75
$ws.Range("M5").Value = -851.4285500000001
$ws.Range("N5").Value = -9050.75
$ws.Range("H92").Value = 818.0909
$ws.Range("I92").Value = 800.6667
$ws.Range("J92").Value = 824.625
$ws.Range("K92").Value = 2402.0001
$ws.Range("L92").Value = 2473.875
$ws.Range("M92").Value = -1154.0001
$ws.Range("N92").Value = -4969.875
$ws.Range("H110").Value = 12205.4
$ws.Range("I110").Value = 4513.5
$ws.Range("J110").Value = 17333.334
$ws.Range("K110").Value = 13540.5
$ws.Range("L110").Value = 52000.00199999999
$ws.Range("M110").Value = -9450.5
$ws.Range("N110").Value = -60180.00199999999
$ws.Range("H111").Value = 12199.4
$ws.Range("I111").Value = 1000
$ws.Range("J111").Value = 14999.25
$ws.Range("K111").Value = 3000
$ws.Range("L111").Value = 44997.75
$ws.Range("M111").Value = 67
$ws.Range("N111").Value = -51131.75
$ws.Range("H112").Value = 4849.512
$ws.Range("I112").Value = 5776
$ws.Range("J112").Value = 4720.8335
$ws.Range("K112").Value = 17328
$ws.Range("L112").Value = 14162.5005
$ws.Range("M112").Value = -16220
$ws.Range("N112").Value = -16378.5005
$ws.Range("H113").Value = 753.15
$ws.Range("I113").Value = 442.5
$ws.Range("J113").Value = 1685.1
$ws.Range("K113").Value = 1327.5
$ws.Range("L113").Value = 5055.299999999999
$ws.Range("M113").Value = 842.5
$ws.Range("N113").Value = -9395.299999999999
$ws.Range("H118").Value = 2032.875
$ws.Range("I118").Value = 2132.8333
$ws.Range("J118").Value = 1999.5555
$ws.Range("K118").Value = 6398.499899999999
$ws.Range("L118").Value = 5998.666499999999
$ws.Range("M118").Value = -5155.499899999999
$ws.Range("N118").Value = -8484.666499999999
$ws.Range("H135").Value = 1530.8846
$ws.Range("I135").Value = 321.14285
$ws.Range("J135").Value = 2942.25
$ws.Range("K135").Value = 2890.28565
$ws.Range("L135").Value = 26480.25
$ws.Range("M135").Value = -355.2856500000003
$ws.Range("N135").Value = -31550.25
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H4").Value = 3500
$ws.Range("I4").Value = 2000
$ws.Range("K4").Value = 2000
$ws.Range("M4").Value = -1888
$ws.Range("H44").Value = 8500
$ws.Range("I44").Value = 8500
$ws.Range("K44").Value = 8500
$ws.Range("M44").Value = -7904
$ws.Range("H102").Value = 3958.8462
$ws.Range("I102").Value = 4072.75
$ws.Range("J102").Value = 2592
$ws.Range("K102").Value = 4072.75
$ws.Range("L102").Value = 2592
$ws.Range("M102").Value = -2450.75
$ws.Range("N102").Value = -5836
$ws.Range("H107").Value = 5342.4
$ws.Range("I107").Value = 6053.5
$ws.Range("J107").Value = 2498
$ws.Range("K107").Value = 6053.5
$ws.Range("L107").Value = 2498
$ws.Range("M107").Value = -4133.5
$ws.Range("N107").Value = -6338
$ws.Range("H108").Value = 50000
$ws.Range("J108").Value = 50000
$ws.Range("L108").Value = 50000
$ws.Range("N108").Value = -57680
$ws.Range("H110").Value = 36368.668
$ws.Range("J110").Value = 36368.668
$ws.Range("L110").Value = 36368.668
$ws.Range("N110").Value = -44548.668
$ws.Range("H112").Value = 45000
$ws.Range("J112").Value = 45000
$ws.Range("L112").Value = 45000
$ws.Range("N112").Value = -47216
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H56").Value = 22838
$ws.Range("I56").Value = 28000
$ws.Range("J56").Value = 21805.6
$ws.Range("K56").Value = 28000
$ws.Range("L56").Value = 21805.6
$ws.Range("M56").Value = -27309
$ws.Range("N56").Value = -23187.6
$ws.Range("H100").Value = 2560.125
$ws.Range("I100").Value = 2650
$ws.Range("J100").Value = 2506.2
$ws.Range("K100").Value = 2650
$ws.Range("L100").Value = 2506.2
$ws.Range("M100").Value = -2109
$ws.Range("N100").Value = -3588.2
